$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.263.05'
$ws.Range('E2').Value = '  +1.54%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.862.40'
$ws.Range('E3').Value = '  +1.32%  '

# Row 4
$ws.Range('E4').Value = '  +1.37%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '312.78'

# Row 6
$ws.Range('E6').Value = '  +1.37%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4800'
$ws.Range('E7').Value = '  +2.02%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3728'
$ws.Range('E8').Value = '  +2.06%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07411'
$ws.Range('E9').Value = '  +3.74%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9377'
$ws.Range('E10').Value = '  +1.86%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.74'
$ws.Range('E11').Value = '  +6.14%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07880'
$ws.Range('E12').Value = '  +2.93%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.871.36'
$ws.Range('E13').Value = '  +3.22%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.435'
$ws.Range('E14').Value = '  +2.87%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.539'
$ws.Range('E15').Value = '  +2.33%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '90.25'
$ws.Range('E16').Value = '  +2.46%  '

# Row 17
$ws.Range('E17').Value = '  +1.41%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008772'
$ws.Range('E18').Value = '  +1.66%  '

# Row 19
$ws.Range('E19').Value = '  +1.38%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.79'
$ws.Range('E20').Value = '  +2.42%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '27.302.79'
$ws.Range('E21').Value = '  +1.59%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.116'
$ws.Range('E22').Value = '  +2.20%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.69'
$ws.Range('E23').Value = '  +0.92%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.953'
$ws.Range('E24').Value = '  +1.77%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.11'
$ws.Range('E25').Value = '  +1.69%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '18.53'
$ws.Range('E26').Value = '  +1.85%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.005'
$ws.Range('E27').Value = '  -0.06%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '115.97'
$ws.Range('E28').Value = '  +1.65%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.997'
$ws.Range('E29').Value = '  +2.54%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.08902'
$ws.Range('E30').Value = '  +0.94%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.346'
$ws.Range('E31').Value = '  +4.23%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.196'
$ws.Range('E32').Value = '  +1.78%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.564'
$ws.Range('E33').Value = '  +2.14%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7436'
$ws.Range('E34').Value = '  -0.30%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.679'
$ws.Range('E35').Value = '  -2.49%  '

# Row 36
$ws.Range('B36').Value = 'VeChain'
$ws.Range('C36').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02048'
$ws.Range('E36').Value = '  +5.48%  '

# Row 37
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.125'
$ws.Range('E37').Value = '  +3.59%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05288'
$ws.Range('E38').Value = '  +1.52%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5369'
$ws.Range('E39').Value = '  +3.37%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '7.124'
$ws.Range('E40').Value = '  +2.41%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.1538'
$ws.Range('E41').Value = '  +1.77%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.374'
$ws.Range('E42').Value = '  +2.80%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '10.60'
$ws.Range('E43').Value = '  +1.48%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.4809'
$ws.Range('E44').Value = '  +2.43%  '

# Row 45
$ws.Range('E45').Value = '  +1.48%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '102.92'
$ws.Range('E46').Value = '  +1.44%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.654'
$ws.Range('E47').Value = '  +3.85%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '66.74'
$ws.Range('E48').Value = '  +2.77%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.06078'
$ws.Range('E49').Value = '  +0.75%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.9011'
$ws.Range('E50').Value = '  +1.74%  '

# Row 51
$ws.Range('E51').Value = '  +1.69%  '
